# Restructure "Feature Key" column (A) from plain numbers to F<n> text
# labels (shared strings), and move the active selection to B16.
#
# Values are written in this specific order so the newly-created shared
# strings land at the same indices as the canonical workbook: starting
# at row 9 and walking upward to row 3, then starting at row 10 and
# walking downward to row 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value  = "F7"
$ws.Range("A8").Value  = "F6"
$ws.Range("A7").Value  = "F5"
$ws.Range("A6").Value  = "F4"
$ws.Range("A5").Value  = "F3"
$ws.Range("A4").Value  = "F2"
$ws.Range("A3").Value  = "F1"
$ws.Range("A10").Value = "F8"
$ws.Range("A11").Value = "F9"
$ws.Range("A12").Value = "F10"
$ws.Range("A13").Value = "F11"
$ws.Range("A14").Value = "F12"
$ws.Range("A15").Value = "F13"
$ws.Range("A16").Value = "F14"
$ws.Range("A17").Value = "F15"

$ws.Range("B16").Select()
